$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -12
$ws.Range("F3").Value = -3
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = 4
$ws.Range("F8").Value = -8
$ws.Range("F16").Value = 5
$ws.Range("F19").Value = -3
$ws.Range("F21").Value = -9
$ws.Range("F23").Value = 5
$ws.Range("F24").Value = -1
$ws.Range("F28").Value = 1
$ws.Range("F30").Value = -1
